$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, row 1
$ws.Range("H1").Value = "Save"

# Match the header formatting used by the existing header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats

# Add the corresponding data value for row 2
$ws.Range("H2").Value = 0
